$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 44003080
$ws.Range("I76").Value = 44003080
$ws.Range("K76").Value = 44003080
$ws.Range("M76").Value = -44002765

$ws.Range("H79").Value = 44003080
$ws.Range("I79").Value = 44003080
$ws.Range("K79").Value = 44003080
$ws.Range("M79").Value = -44001988

$ws.Range("H86").Value = 83336050
$ws.Range("I86").Value = 882.6667
$ws.Range("J86").Value = 166671220
$ws.Range("K86").Value = 882.6667
$ws.Range("L86").Value = 166671220
$ws.Range("M86").Value = 240.3333
$ws.Range("N86").Value = -166673466

$ws.Range("H89").Value = 83336050
$ws.Range("I89").Value = 882.6667
$ws.Range("J89").Value = 166671220
$ws.Range("K89").Value = 4413.3335
$ws.Range("L89").Value = 833356100
$ws.Range("M89").Value = 1202.6665
$ws.Range("N89").Value = -833367332

$ws.Range("H113").Value = 2759
$ws.Range("I113").Value = 2759
$ws.Range("K113").Value = 2759
$ws.Range("M113").Value = 495

$ws.Range("H116").Value = 4528.5713
$ws.Range("I116").Value = 4450
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 4450
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = -1008
$ws.Range("N116").Value = -11884

$ws.Range("H129").Value = 1278200.4
$ws.Range("I129").Value = 457.6
$ws.Range("J129").Value = 1544396.8
$ws.Range("K129").Value = 1372.8
$ws.Range("L129").Value = 4633190.4
$ws.Range("M129").Value = 3627.2
$ws.Range("N129").Value = -4643190.4

$ws.Range("H132").Value = 2112.9038
$ws.Range("I132").Value = 1997.2174
$ws.Range("K132").Value = 5991.6522
$ws.Range("M132").Value = -3461.6522

$ws.Range("H138").Value = 3651.48
$ws.Range("I138").Value = 827.1
$ws.Range("J138").Value = 4357.575
$ws.Range("K138").Value = 2481.3
$ws.Range("L138").Value = 13072.725
$ws.Range("M138").Value = 2658.7
$ws.Range("N138").Value = -23352.725

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 907.9459000000001
$ws.Range("I61").Value = 752.5925999999999
$ws.Range("J61").Value = 1327.4
$ws.Range("K61").Value = 752.5925999999999
$ws.Range("L61").Value = 1327.4
$ws.Range("M61").Value = -540.5925999999999
$ws.Range("N61").Value = -1751.4

$ws.Range("H74").Value = 1048.25
$ws.Range("I74").Value = 979.4706
$ws.Range("J74").Value = 1438
$ws.Range("K74").Value = 979.4706
$ws.Range("L74").Value = 1438
$ws.Range("M74").Value = -105.4706
$ws.Range("N74").Value = -3186

$ws.Range("H77").Value = 1048.25
$ws.Range("I77").Value = 979.4706
$ws.Range("J77").Value = 1438
$ws.Range("K77").Value = 4897.353
$ws.Range("L77").Value = 7190
$ws.Range("M77").Value = -529.3530000000001
$ws.Range("N77").Value = -15926

$ws.Range("H132").Value = 1169.9714
$ws.Range("I132").Value = 808.9583
$ws.Range("J132").Value = 1957.6364
$ws.Range("K132").Value = 2426.8749
$ws.Range("L132").Value = 5872.9092
$ws.Range("M132").Value = 103.1251000000002
$ws.Range("N132").Value = -10932.9092

$ws.Range("H136").Value = 907.9459000000001
$ws.Range("I136").Value = 752.5925999999999
$ws.Range("J136").Value = 1327.4
$ws.Range("K136").Value = 2257.7778
$ws.Range("L136").Value = 3982.2
$ws.Range("M136").Value = 292.2222000000002
$ws.Range("N136").Value = -9082.200000000001

$ws.Range("H139").Value = 47626.668
$ws.Range("J139").Value = 47626.668
$ws.Range("L139").Value = 47626.668
$ws.Range("N139").Value = -57906.668

$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("M140").ClearContents()
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 101714.29
$ws.Range("I20").Value = 151500
$ws.Range("J20").Value = 35333.332
$ws.Range("K20").Value = 151500
$ws.Range("L20").Value = 35333.332
$ws.Range("M20").Value = -151253
$ws.Range("N20").Value = -35827.332

$ws.Range("H22").Value = 290
$ws.Range("I22").Value = 290
$ws.Range("K22").Value = 290
$ws.Range("M22").Value = -117

$ws.Range("H86").Value = 2053.8965
$ws.Range("I86").Value = 1870.2609
$ws.Range("K86").Value = 1870.2609
$ws.Range("M86").Value = -747.2609

$ws.Range("H89").Value = 2053.8965
$ws.Range("I89").Value = 1870.2609
$ws.Range("K89").Value = 9351.3045
$ws.Range("M89").Value = -3735.3045

$ws.Range("H105").Value = 7943.4
$ws.Range("I105").Value = 9226.25
$ws.Range("J105").Value = 6477.2856
$ws.Range("K105").Value = 9226.25
$ws.Range("L105").Value = 6477.2856
$ws.Range("M105").Value = -7479.25
$ws.Range("N105").Value = -9971.285599999999

$ws.Range("H134").Value = 82393.12
$ws.Range("I134").Value = 2737.5
$ws.Range("J134").Value = 224003.11
$ws.Range("K134").Value = 8212.5
$ws.Range("L134").Value = 672009.33
$ws.Range("M134").Value = -5677.5
$ws.Range("N134").Value = -677079.33

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1793.6086
$ws.Range("I132").Value = 1071.25
$ws.Range("K132").Value = 3213.75
$ws.Range("M132").Value = -683.75

$ws.Range("H134").Value = 2520.946
$ws.Range("I134").Value = 1759.3572
$ws.Range("J134").Value = 4890.3335
$ws.Range("K134").Value = 5278.071599999999
$ws.Range("L134").Value = 14671.0005
$ws.Range("M134").Value = -2743.071599999999
$ws.Range("N134").Value = -19741.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 1023.5
$ws.Range("I6").Value = 149.11111
$ws.Range("J6").Value = 3646.6667
$ws.Range("K6").Value = 447.33333
$ws.Range("L6").Value = 10940.0001
$ws.Range("M6").Value = -334.33333
$ws.Range("N6").Value = -11166.0001

$ws.Range("H131").Value = 817.86
$ws.Range("J131").Value = 825.7755
$ws.Range("L131").Value = 2477.3265
$ws.Range("N131").Value = -12557.3265

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4361.6895
$ws.Range("I70").Value = 4087.5
$ws.Range("K70").Value = 4087.5
$ws.Range("M70").Value = -3817.5

$ws.Range("H73").Value = 4361.6895
$ws.Range("I73").Value = 4087.5
$ws.Range("K73").Value = 4087.5
$ws.Range("M73").Value = -3151.5

$ws.Range("H132").Value = 3629.4736
$ws.Range("I132").Value = 3777.0908
$ws.Range("J132").Value = 3426.5
$ws.Range("K132").Value = 11331.2724
$ws.Range("L132").Value = 10279.5
$ws.Range("M132").Value = -8801.2724
$ws.Range("N132").Value = -15339.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5904.2666
$ws.Range("I122").Value = 6436.72
$ws.Range("J122").Value = 3242
$ws.Range("K122").Value = 19310.16
$ws.Range("L122").Value = 9726
$ws.Range("M122").Value = -16860.16
$ws.Range("N122").Value = -14626
